$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "24.03.2023 15:41 (CET)"
$ws.Range("C5").Value = "https://gitlab.intra.infineon.com/semantic-web-projects/digital-reference/order_management/-/commit/fb965955f5d3cc29931426bf15edfcd279305895"
$ws.Range("D5").Value = "b95ad415b600b67a1ceb8669c0838ed63b287ee61121c14e73fa33c6d73200af"
